$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 706.5  # H76: 0 -> 706.5
$ws.Cells.Item(76, 9).Value = 413  # I76: 0 -> 413
$ws.Cells.Item(76, 10).Value = 1000  # J76: 0 -> 1000
$ws.Cells.Item(76, 11).Value = 413  # K76: 0 -> 413
$ws.Cells.Item(76, 12).Value = 1000  # L76: 0 -> 1000
$ws.Cells.Item(76, 13).Value = -98  # M76: None -> -98
$ws.Cells.Item(76, 14).Value = -1630  # N76: None -> -1630

$ws.Cells.Item(79, 8).Value = 706.5  # H79: 0 -> 706.5
$ws.Cells.Item(79, 9).Value = 413  # I79: 0 -> 413
$ws.Cells.Item(79, 10).Value = 1000  # J79: 0 -> 1000
$ws.Cells.Item(79, 11).Value = 413  # K79: 0 -> 413
$ws.Cells.Item(79, 12).Value = 1000  # L79: 0 -> 1000
$ws.Cells.Item(79, 13).Value = 679  # M79: None -> 679
$ws.Cells.Item(79, 14).Value = -3184  # N79: None -> -3184

$ws.Cells.Item(80, 8).Value = 763.7778  # H80: 798.1429000000001 -> 763.7778
$ws.Cells.Item(80, 9).Value = 645.6667  # I80: 646.75 -> 645.6667
$ws.Cells.Item(80, 11).Value = 1937.0001  # K80: 1940.25 -> 1937.0001
$ws.Cells.Item(80, 13).Value = -939.0001  # M80: -942.25 -> -939.0001

$ws.Cells.Item(83, 8).Value = 763.7778  # H83: 798.1429000000001 -> 763.7778
$ws.Cells.Item(83, 9).Value = 645.6667  # I83: 646.75 -> 645.6667
$ws.Cells.Item(83, 11).Value = 5811.0003  # K83: 5820.75 -> 5811.0003
$ws.Cells.Item(83, 13).Value = -819.0002999999997  # M83: -828.75 -> -819.0002999999997

$ws.Cells.Item(87, 8).Value = 94784.336  # H87: 95338.25 -> 94784.336
$ws.Cells.Item(87, 10).Value = 94784.336  # J87: 95338.25 -> 94784.336
$ws.Cells.Item(87, 12).Value = 94784.336  # L87: 95338.25 -> 94784.336
$ws.Cells.Item(87, 14).Value = -97280.336  # N87: -97834.25 -> -97280.336

$ws.Cells.Item(88, 8).Value = 1282.25  # H88: 1665.7778 -> 1282.25
$ws.Cells.Item(88, 9).Value = 1419.5  # I88: 1799.2 -> 1419.5
$ws.Cells.Item(88, 10).Value = 1145  # J88: 1499 -> 1145
$ws.Cells.Item(88, 11).Value = 1419.5  # K88: 1799.2 -> 1419.5
$ws.Cells.Item(88, 12).Value = 1145  # L88: 1499 -> 1145
$ws.Cells.Item(88, 13).Value = -1013.5  # M88: -1393.2 -> -1013.5
$ws.Cells.Item(88, 14).Value = -1957  # N88: -2311 -> -1957

$ws.Cells.Item(90, 8).Value = 94784.336  # H90: 95338.25 -> 94784.336
$ws.Cells.Item(90, 10).Value = 94784.336  # J90: 95338.25 -> 94784.336
$ws.Cells.Item(90, 12).Value = 284353.008  # L90: 286014.75 -> 284353.008
$ws.Cells.Item(90, 14).Value = -296833.008  # N90: -298494.75 -> -296833.008

$ws.Cells.Item(91, 8).Value = 1282.25  # H91: 1665.7778 -> 1282.25
$ws.Cells.Item(91, 9).Value = 1419.5  # I91: 1799.2 -> 1419.5
$ws.Cells.Item(91, 10).Value = 1145  # J91: 1499 -> 1145
$ws.Cells.Item(91, 11).Value = 1419.5  # K91: 1799.2 -> 1419.5
$ws.Cells.Item(91, 12).Value = 1145  # L91: 1499 -> 1145
$ws.Cells.Item(91, 13).Value = -15.5  # M91: -395.2 -> -15.5
$ws.Cells.Item(91, 14).Value = -3953  # N91: -4307 -> -3953

$ws.Cells.Item(113, 8).Value = 7601  # H113: 8251.25 -> 7601
$ws.Cells.Item(113, 9).Value = 7601  # I113: 8251.25 -> 7601
$ws.Cells.Item(113, 11).Value = 7601  # K113: 8251.25 -> 7601
$ws.Cells.Item(113, 13).Value = -4347  # M113: -4997.25 -> -4347

$ws.Cells.Item(127, 8).Value = 0  # H127: 2200 -> 0
$ws.Cells.Item(127, 9).Value = 0  # I127: 2200 -> 0
$ws.Cells.Item(127, 11).Value = 0  # K127: 6600 -> 0
$ws.Cells.Item(127, 13).ClearContents()  # M127: -1640 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(59, 8).Value = 10059  # H59: 0 -> 10059
$ws.Cells.Item(59, 10).Value = 10059  # J59: 0 -> 10059
$ws.Cells.Item(59, 12).Value = 10059  # L59: 0 -> 10059
$ws.Cells.Item(59, 14).Value = -11667  # N59: None -> -11667

$ws.Cells.Item(88, 8).Value = 1384.7  # H88: 1718.8572 -> 1384.7
$ws.Cells.Item(88, 9).Value = 1506  # I88: 1506.3334 -> 1506
$ws.Cells.Item(88, 10).Value = 1303.8334  # J88: 1878.25 -> 1303.8334
$ws.Cells.Item(88, 11).Value = 1506  # K88: 1506.3334 -> 1506
$ws.Cells.Item(88, 12).Value = 1303.8334  # L88: 1878.25 -> 1303.8334
$ws.Cells.Item(88, 13).Value = -1100  # M88: -1100.3334 -> -1100
$ws.Cells.Item(88, 14).Value = -2115.8334  # N88: -2690.25 -> -2115.8334

$ws.Cells.Item(91, 8).Value = 1384.7  # H91: 1718.8572 -> 1384.7
$ws.Cells.Item(91, 9).Value = 1506  # I91: 1506.3334 -> 1506
$ws.Cells.Item(91, 10).Value = 1303.8334  # J91: 1878.25 -> 1303.8334
$ws.Cells.Item(91, 11).Value = 1506  # K91: 1506.3334 -> 1506
$ws.Cells.Item(91, 12).Value = 1303.8334  # L91: 1878.25 -> 1303.8334
$ws.Cells.Item(91, 13).Value = -102  # M91: -102.3334 -> -102
$ws.Cells.Item(91, 14).Value = -4111.8334  # N91: -4686.25 -> -4111.8334

$ws.Cells.Item(122, 8).Value = 4502.5  # H122: 5000 -> 4502.5
$ws.Cells.Item(122, 9).Value = 4502.5  # I122: 5000 -> 4502.5
$ws.Cells.Item(122, 11).Value = 13507.5  # K122: 15000 -> 13507.5
$ws.Cells.Item(122, 13).Value = -11057.5  # M122: -12550 -> -11057.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1821  # H105: 1859 -> 1821
$ws.Cells.Item(105, 9).Value = 1869  # I105: 1905.2 -> 1869
$ws.Cells.Item(105, 10).Value = 1725  # J105: 1743.5 -> 1725
$ws.Cells.Item(105, 11).Value = 1869  # K105: 1905.2 -> 1869
$ws.Cells.Item(105, 12).Value = 1725  # L105: 1743.5 -> 1725
$ws.Cells.Item(105, 13).Value = -122  # M105: -158.2 -> -122
$ws.Cells.Item(105, 14).Value = -5219  # N105: -5237.5 -> -5219

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9599.143  # H31: 7563.6 -> 9599.143
$ws.Cells.Item(31, 9).Value = 6500  # I31: 3064 -> 6500
$ws.Cells.Item(31, 10).Value = 9837.538  # J31: 9492 -> 9837.538
$ws.Cells.Item(31, 11).Value = 6500  # K31: 3064 -> 6500
$ws.Cells.Item(31, 12).Value = 9837.538  # L31: 9492 -> 9837.538
$ws.Cells.Item(31, 13).Value = -6205  # M31: -2769 -> -6205
$ws.Cells.Item(31, 14).Value = -10427.538  # N31: -10082 -> -10427.538

$ws.Cells.Item(34, 8).Value = 9599.143  # H34: 7563.6 -> 9599.143
$ws.Cells.Item(34, 9).Value = 6500  # I34: 3064 -> 6500
$ws.Cells.Item(34, 10).Value = 9837.538  # J34: 9492 -> 9837.538
$ws.Cells.Item(34, 11).Value = 6500  # K34: 3064 -> 6500
$ws.Cells.Item(34, 12).Value = 9837.538  # L34: 9492 -> 9837.538
$ws.Cells.Item(34, 13).Value = -6298  # M34: -2862 -> -6298
$ws.Cells.Item(34, 14).Value = -10241.538  # N34: -9896 -> -10241.538

$ws.Cells.Item(36, 8).Value = 0  # H36: 1366.6666 -> 0
$ws.Cells.Item(36, 9).Value = 0  # I36: 1366.6666 -> 0
$ws.Cells.Item(36, 11).Value = 0  # K36: 1366.6666 -> 0
$ws.Cells.Item(36, 13).ClearContents()  # M36: -978.6666 -> (removed)

$ws.Cells.Item(40, 8).Value = 0  # H40: 1366.6666 -> 0
$ws.Cells.Item(40, 9).Value = 0  # I40: 1366.6666 -> 0
$ws.Cells.Item(40, 11).Value = 0  # K40: 1366.6666 -> 0
$ws.Cells.Item(40, 13).ClearContents()  # M40: -1206.6666 -> (removed)

$ws.Cells.Item(58, 8).Value = 8496.25  # H58: 8665 -> 8496.25
$ws.Cells.Item(58, 9).Value = 7995  # I58: 8000 -> 7995
$ws.Cells.Item(58, 11).Value = 7995  # K58: 8000 -> 7995
$ws.Cells.Item(58, 13).Value = -7792  # M58: -7797 -> -7792

$ws.Cells.Item(68, 8).Value = 54747.8  # H68: 46794 -> 54747.8
$ws.Cells.Item(68, 9).Value = 27999  # I68: 27999.666 -> 27999
$ws.Cells.Item(68, 10).Value = 61435  # J68: 65588.336 -> 61435
$ws.Cells.Item(68, 11).Value = 27999  # K68: 27999.666 -> 27999
$ws.Cells.Item(68, 12).Value = 61435  # L68: 65588.336 -> 61435
$ws.Cells.Item(68, 13).Value = -27250  # M68: -27250.666 -> -27250
$ws.Cells.Item(68, 14).Value = -62933  # N68: -67086.336 -> -62933

$ws.Cells.Item(71, 8).Value = 54747.8  # H71: 46794 -> 54747.8
$ws.Cells.Item(71, 9).Value = 27999  # I71: 27999.666 -> 27999
$ws.Cells.Item(71, 10).Value = 61435  # J71: 65588.336 -> 61435
$ws.Cells.Item(71, 11).Value = 83997  # K71: 83998.99800000001 -> 83997
$ws.Cells.Item(71, 12).Value = 184305  # L71: 196765.008 -> 184305
$ws.Cells.Item(71, 13).Value = -80253  # M71: -80254.99800000001 -> -80253
$ws.Cells.Item(71, 14).Value = -191793  # N71: -204253.008 -> -191793

$ws.Cells.Item(99, 8).Value = 1930.0714  # H99: 1936.3077 -> 1930.0714
$ws.Cells.Item(99, 9).Value = 1880.25  # I99: 1890.6666 -> 1880.25
$ws.Cells.Item(99, 11).Value = 1880.25  # K99: 1890.6666 -> 1880.25
$ws.Cells.Item(99, 13).Value = -382.25  # M99: -392.6666 -> -382.25

$ws.Cells.Item(126, 8).Value = 1930.0714  # H126: 1936.3077 -> 1930.0714
$ws.Cells.Item(126, 9).Value = 1880.25  # I126: 1890.6666 -> 1880.25
$ws.Cells.Item(126, 11).Value = 5640.75  # K126: 5671.9998 -> 5640.75
$ws.Cells.Item(126, 13).Value = -3170.75  # M126: -3201.9998 -> -3170.75

$ws.Cells.Item(136, 8).Value = 8496.25  # H136: 8665 -> 8496.25
$ws.Cells.Item(136, 9).Value = 7995  # I136: 8000 -> 7995
$ws.Cells.Item(136, 11).Value = 23985  # K136: 24000 -> 23985
$ws.Cells.Item(136, 13).Value = -21435  # M136: -21450 -> -21435

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 161.81818  # H40: 216 -> 161.81818
$ws.Cells.Item(40, 9).Value = 25.714285  # I40: 23.5 -> 25.714285
$ws.Cells.Item(40, 10).Value = 400  # J40: 370 -> 400
$ws.Cells.Item(40, 11).Value = 102.85714  # K40: 94 -> 102.85714
$ws.Cells.Item(40, 12).Value = 1600  # L40: 1480 -> 1600
$ws.Cells.Item(40, 13).Value = -33.85714  # M40: -25 -> -33.85714
$ws.Cells.Item(40, 14).Value = -1738  # N40: -1618 -> -1738

$ws.Cells.Item(132, 8).Value = 2891.923  # H132: 2907.1538 -> 2891.923
$ws.Cells.Item(132, 9).Value = 2463.182  # I132: 2481.182 -> 2463.182
$ws.Cells.Item(132, 11).Value = 22168.638  # K132: 22330.638 -> 22168.638
$ws.Cells.Item(132, 13).Value = -19638.638  # M132: -19800.638 -> -19638.638

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 12128296  # H11: 7687784 -> 12128296
$ws.Cells.Item(11, 9).Value = 8646776  # I11: 4871689 -> 8646776
$ws.Cells.Item(11, 11).Value = 8646776  # K11: 4871689 -> 8646776
$ws.Cells.Item(11, 13).Value = -8646637  # M11: -4871550 -> -8646637

$ws.Cells.Item(24, 8).Value = 44444  # H24: 65007 -> 44444
$ws.Cells.Item(24, 10).Value = 44444  # J24: 65007 -> 44444
$ws.Cells.Item(24, 12).Value = 44444  # L24: 65007 -> 44444
$ws.Cells.Item(24, 14).Value = -44790  # N24: -65353 -> -44790

$ws.Cells.Item(39, 8).Value = 0  # H39: 22222 -> 0
$ws.Cells.Item(39, 10).Value = 0  # J39: 22222 -> 0
$ws.Cells.Item(39, 12).Value = 0  # L39: 22222 -> 0
$ws.Cells.Item(39, 14).ClearContents()  # N39: -23286 -> (removed)

$ws.Cells.Item(132, 8).Value = 1025  # H132: 1041.4445 -> 1025
$ws.Cells.Item(132, 10).Value = 0  # J132: 1099 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 3297 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -8357 -> (removed)

$ws.Cells.Item(140, 8).Value = 111505.5  # H140: 101391.86 -> 111505.5
$ws.Cells.Item(140, 9).Value = 279697  # I140: 160203 -> 279697
$ws.Cells.Item(140, 10).Value = 77867.2  # J140: 77867.39999999999 -> 77867.2
$ws.Cells.Item(140, 11).Value = 279697  # K140: 160203 -> 279697
$ws.Cells.Item(140, 12).Value = 77867.2  # L140: 77867.39999999999 -> 77867.2
$ws.Cells.Item(140, 13).Value = -274517  # M140: -155023 -> -274517
$ws.Cells.Item(140, 14).Value = -88227.2  # N140: -88227.39999999999 -> -88227.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2199.6667  # H22: 1000 -> 2199.6667
$ws.Cells.Item(22, 9).Value = 1549.5  # I22: 0 -> 1549.5
$ws.Cells.Item(22, 10).Value = 3500  # J22: 1000 -> 3500
$ws.Cells.Item(22, 11).Value = 1549.5  # K22: 0 -> 1549.5
$ws.Cells.Item(22, 12).Value = 3500  # L22: 1000 -> 3500
$ws.Cells.Item(22, 13).Value = -1254.5  # M22: None -> -1254.5
$ws.Cells.Item(22, 14).Value = -4090  # N22: -1590 -> -4090

$ws.Cells.Item(27, 8).Value = 2199.6667  # H27: 1000 -> 2199.6667
$ws.Cells.Item(27, 9).Value = 1549.5  # I27: 0 -> 1549.5
$ws.Cells.Item(27, 10).Value = 3500  # J27: 1000 -> 3500
$ws.Cells.Item(27, 11).Value = 1549.5  # K27: 0 -> 1549.5
$ws.Cells.Item(27, 12).Value = 3500  # L27: 1000 -> 3500
$ws.Cells.Item(27, 13).Value = -1442.5  # M27: None -> -1442.5
$ws.Cells.Item(27, 14).Value = -3714  # N27: -1214 -> -3714

$ws.Cells.Item(44, 8).Value = 20000  # H44: 7500 -> 20000
$ws.Cells.Item(44, 9).Value = 0  # I44: 5000 -> 0
$ws.Cells.Item(44, 10).Value = 20000  # J44: 10000 -> 20000
$ws.Cells.Item(44, 11).Value = 0  # K44: 5000 -> 0
$ws.Cells.Item(44, 12).Value = 20000  # L44: 10000 -> 20000
$ws.Cells.Item(44, 13).ClearContents()  # M44: -4544 -> (removed)
$ws.Cells.Item(44, 14).Value = -20912  # N44: -10912 -> -20912

$ws.Cells.Item(61, 8).Value = 3778.2856  # H61: 4324.8335 -> 3778.2856
$ws.Cells.Item(61, 9).Value = 612  # I61: 649.6667 -> 612
$ws.Cells.Item(61, 11).Value = 612  # K61: 649.6667 -> 612
$ws.Cells.Item(61, 13).Value = -410  # M61: -447.6667 -> -410

$ws.Cells.Item(113, 8).Value = 3778.2856  # H113: 4324.8335 -> 3778.2856
$ws.Cells.Item(113, 9).Value = 612  # I113: 649.6667 -> 612
$ws.Cells.Item(113, 11).Value = 612  # K113: 649.6667 -> 612
$ws.Cells.Item(113, 13).Value = 1558  # M113: 1520.3333 -> 1558

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 428.8  # H81: 448.8 -> 428.8
$ws.Cells.Item(81, 9).Value = 428.8  # I81: 373.5 -> 428.8
$ws.Cells.Item(81, 10).Value = 0  # J81: 750 -> 0
$ws.Cells.Item(81, 11).Value = 857.6  # K81: 747 -> 857.6
$ws.Cells.Item(81, 12).Value = 0  # L81: 1500 -> 0
$ws.Cells.Item(81, 13).Value = 203.4  # M81: 314 -> 203.4
$ws.Cells.Item(81, 14).ClearContents()  # N81: -3622 -> (removed)

$ws.Cells.Item(84, 8).Value = 428.8  # H84: 448.8 -> 428.8
$ws.Cells.Item(84, 9).Value = 428.8  # I84: 373.5 -> 428.8
$ws.Cells.Item(84, 10).Value = 0  # J84: 750 -> 0
$ws.Cells.Item(84, 11).Value = 4288  # K84: 3735 -> 4288
$ws.Cells.Item(84, 12).Value = 0  # L84: 7500 -> 0
$ws.Cells.Item(84, 13).Value = 1016  # M84: 1569 -> 1016
$ws.Cells.Item(84, 14).ClearContents()  # N84: -18108 -> (removed)

$ws.Cells.Item(132, 8).Value = 1495  # H132: 1550.2222 -> 1495
$ws.Cells.Item(132, 9).Value = 1220.7142  # I132: 1257.8334 -> 1220.7142
$ws.Cells.Item(132, 11).Value = 3662.1426  # K132: 3773.5002 -> 3662.1426
$ws.Cells.Item(132, 13).Value = -1132.1426  # M132: -1243.5002 -> -1132.1426
